# Applies the dac-val.xlsx update: new Param2 column, refreshed rows,
# renamed/reordered classes, and trimmed row count (26 -> 21 data+header rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Param2" header column (D), formatted like the other header cells ---
$ws.Range('C1').Copy() | Out-Null
$ws.Range('D1').PasteSpecial(-4122) | Out-Null
$ws.Range('D1').Value = 'Param2'

# --- Rewrite all data rows (2-21) with the updated dataset ---
$ws.Range('A2').Value = 'zelda--param1-00.98--1-08 Rune Extraction.dac'
$ws.Range('B2').Value = 'zelda'
$ws.Range('C2').Value = 0.98
$ws.Range('D2').Value = '08 Rune Extraction'

$ws.Range('A3').Value = 'zelda--param1-00.37--5-13 Windblight Ganon Appears.dac'
$ws.Range('B3').Value = 'zelda'
$ws.Range('C3').Value = 0.37
$ws.Range('D3').Value = '13 Windblight Ganon Appears'

$ws.Range('A4').Value = 'zelda--param1-00.00.dac'
$ws.Range('B4').Value = 'zelda'
$ws.Range('C4').Value = 0
$ws.Range('D4').Value = ''

$ws.Range('A5').Value = 'zelda--param1-00.10--1-29 Impa''s Theme.dac'
$ws.Range('B5').Value = 'zelda'
$ws.Range('C5').Value = 0.1
$ws.Range('D5').Value = '29 Impa''s Theme'

$ws.Range('A6').Value = 'zelda--param1-00.71--3-16 Riju''s Thunder Barrier.dac'
$ws.Range('B6').Value = 'zelda'
$ws.Range('C6').Value = 0.71
$ws.Range('D6').Value = '16 Riju''s Thunder Barrier'

$ws.Range('A7').Value = 'fusion--param1-00.28.dac'
$ws.Range('B7').Value = 'fusion'
$ws.Range('C7').Value = 0.28
$ws.Range('D7').Value = ''

$ws.Range('A8').Value = 'fusion--param1-00.03.dac'
$ws.Range('B8').Value = 'fusion'
$ws.Range('C8').Value = 0.03
$ws.Range('D8').Value = ''

$ws.Range('A9').Value = 'fusion--param1-00.94--02. Lament.dac'
$ws.Range('B9').Value = 'fusion'
$ws.Range('C9').Value = 0.94
$ws.Range('D9').Value = ''

$ws.Range('A10').Value = 'fusion--param1-00.08--05. Rising.dac'
$ws.Range('B10').Value = 'fusion'
$ws.Range('C10').Value = 0.08
$ws.Range('D10').Value = ''

$ws.Range('A11').Value = 'fusion--param1-00.59.dac'
$ws.Range('B11').Value = 'fusion'
$ws.Range('C11').Value = 0.59
$ws.Range('D11').Value = ''

$ws.Range('A12').Value = '8bit--param1-00.37.dac'
$ws.Range('B12').Value = '8bit'
$ws.Range('C12').Value = 0.37
$ws.Range('D12').Value = ''

$ws.Range('A13').Value = '8bit--param1-00.28--30 The King''s Curse.dac'
$ws.Range('B13').Value = '8bit'
$ws.Range('C13').Value = 0.28
$ws.Range('D13').Value = ''

$ws.Range('A14').Value = '8bit--param1-00.82.dac'
$ws.Range('B14').Value = '8bit'
$ws.Range('C14').Value = 0.82
$ws.Range('D14').Value = ''

$ws.Range('A15').Value = '8bit--param1-00.52.dac'
$ws.Range('B15').Value = '8bit'
$ws.Range('C15').Value = 0.52
$ws.Range('D15').Value = ''

$ws.Range('A16').Value = '8bit--param1-00.71--24 The Ice Queen.dac'
$ws.Range('B16').Value = '8bit'
$ws.Range('C16').Value = 0.71
$ws.Range('D16').Value = ''

$ws.Range('A17').Value = 'rock--param1-00.30.dac'
$ws.Range('B17').Value = 'rock'
$ws.Range('C17').Value = 0.3
$ws.Range('D17').Value = ''

$ws.Range('A18').Value = 'rock--param1-00.53.dac'
$ws.Range('B18').Value = 'rock'
$ws.Range('C18').Value = 0.53
$ws.Range('D18').Value = ''

$ws.Range('A19').Value = 'rock--param1-00.39.dac'
$ws.Range('B19').Value = 'rock'
$ws.Range('C19').Value = 0.39
$ws.Range('D19').Value = ''

$ws.Range('A20').Value = 'rock--param1-00.57.dac'
$ws.Range('B20').Value = 'rock'
$ws.Range('C20').Value = 0.57
$ws.Range('D20').Value = ''

$ws.Range('A21').Value = 'rock--param1-00.10.dac'
$ws.Range('B21').Value = 'rock'
$ws.Range('C21').Value = 0.1
$ws.Range('D21').Value = ''

# --- Remove the now-unused trailing rows (old rows 22-26, the "duduk" class) ---
$ws.Range('A22:D26').EntireRow.Delete() | Out-Null

$ur = $ws.UsedRange
Write-Host ("UsedRange=" + $ur.Address())
